# Revert "tamoc-119: use hash index (#5033)"
#
# This removes the "visibilityStatus" column that had been added to the
# Patient / Allergy / Diagnosis reference-data export sheets, and restores
# the Patient sheet's sample row to its pre-change values.

$wb = $excel.ActiveWorkbook

# --- Sheet "Patient": drop column K (visibilityStatus); L/M shift left ---
$wsPatient = $wb.Worksheets.Item("Patient")
$wsPatient.Columns("K").Delete()

# Restore the original sample data row (row 2)
$wsPatient.Range("A2").Value = "3cf6cc61-e682-41ba-aa96-7a46ebef428e"
$wsPatient.Range("B2").Value = "NKIX135928"
$wsPatient.Range("C2").Value = "Helena"
$wsPatient.Range("E2").Value = "Gabbrielli"
$wsPatient.Range("F2").Value = "Inoue"
$wsPatient.Range("G2").Value = 38106.5

# Keep the "number stored as text" ignored-error marker over the new extent
$wsPatient.Range("A1:L2").IgnoredErrors.NumberAsText = $true

# --- Sheet "Allergy": drop column D (visibilityStatus) ---
$wsAllergy = $wb.Worksheets.Item("Allergy")
$wsAllergy.Columns("D").Delete()
$wsAllergy.Range("A1:C3").IgnoredErrors.NumberAsText = $true

# --- Sheet "Diagnosis": drop column D (visibilityStatus) ---
$wsDiagnosis = $wb.Worksheets.Item("Diagnosis")
$wsDiagnosis.Columns("D").Delete()
$wsDiagnosis.Range("A1:C3").IgnoredErrors.NumberAsText = $true
